$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$co = $ws.ChartObjects().Add(100, 20, 300, 200)
$chart = $co.Chart
$chart.ChartType = 5
$chart.SetSourceData($ws.Range("D10:E14"))
$chart.HasDataLabels = $true
Write-Host "chart HasDataLabels set"
